$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (source row 85)
$ws.Range("D2").Value = 45113
$ws.Range("J2").Value = 410
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 18878
$ws.Range("N2").Value = '$/malla 15 kilos'
$ws.Range("P2").Value = 1259
$ws.Range("Q2").Value = 15

# Row 3 (source row 55)
$ws.Range("D3").Value = 44832
$ws.Range("J3").Value = 160
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("O3").Value = 'Hijuelas'
$ws.Range("P3").Value = 1000

# Row 4 (source row 56)
$ws.Range("D4").Value = 44832
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 19000
$ws.Range("M4").Value = 18080
$ws.Range("N4").Value = '$/malla 15 kilos'
$ws.Range("P4").Value = 1205
$ws.Range("Q4").Value = 15

# Row 5 (source row 100)
$ws.Range("D5").Value = 45141
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17545
$ws.Range("P5").Value = 1170

# Row 6 (source row 67)
$ws.Range("D6").Value = 45099
$ws.Range("J6").Value = 140
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 17000
$ws.Range("P6").Value = 1133

# Row 7 (source row 54)
$ws.Range("D7").Value = 44396
$ws.Range("J7").Value = 130
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 22000
$ws.Range("P7").Value = 1467

# Row 8 (source row 59)
$ws.Range("D8").Value = 44483
$ws.Range("J8").Value = 220
$ws.Range("M8").Value = 18909
$ws.Range("P8").Value = 1261

# Row 9 (source row 37)
$ws.Range("D9").Value = 45054
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 20800
$ws.Range("P9").Value = 1387

# Row 10 (source row 80)
$ws.Range("D10").Value = 45062
$ws.Range("J10").Value = 300
$ws.Range("M10").Value = 18800
$ws.Range("P10").Value = 1253

# Row 11 (source row 7)
$ws.Range("D11").Value = 45079
$ws.Range("J11").Value = 270
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19037
$ws.Range("P11").Value = 1269

# Row 12 (source row 89)
$ws.Range("D12").Value = 44812
$ws.Range("J12").Value = 410
$ws.Range("M12").Value = 17488
$ws.Range("P12").Value = 1166

# Row 13 (source row 115)
$ws.Range("D13").Value = 44754
$ws.Range("J13").Value = 300
$ws.Range("L13").Value = 19000
$ws.Range("M13").Value = 18133
$ws.Range("P13").Value = 1209

# Row 14 (source row 84)
$ws.Range("D14").Value = 45072
$ws.Range("J14").Value = 200
$ws.Range("M14").Value = 18800
$ws.Range("P14").Value = 1253

# Row 15 (source row 77)
$ws.Range("D15").Value = 45106
$ws.Range("J15").Value = 580
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12448
$ws.Range("P15").Value = 830

# Row 16 (source row 111)
$ws.Range("D16").Value = 44792
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 16000
$ws.Range("M16").Value = 16800
$ws.Range("P16").Value = 1120

# Row 17 (source row 94)
$ws.Range("D17").Value = 44763
$ws.Range("J17").Value = 350
$ws.Range("M17").Value = 17571
$ws.Range("P17").Value = 1171

# Row 18 (source row 34)
$ws.Range("D18").Value = 44839
$ws.Range("J18").Value = 280
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 19000
$ws.Range("M18").Value = 17857
$ws.Range("P18").Value = 1190

# Row 19 (source row 14)
$ws.Range("D19").Value = 45132
$ws.Range("J19").Value = 290
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 19172
$ws.Range("P19").Value = 1278

# Row 20 (source row 21)
$ws.Range("D20").Value = 45155
$ws.Range("J20").Value = 330
$ws.Range("L20").Value = 19000
$ws.Range("M20").Value = 18394
$ws.Range("P20").Value = 1226

# Row 21 (source row 88)
$ws.Range("D21").Value = 45119
$ws.Range("J21").Value = 157
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17236
$ws.Range("P21").Value = 1149

# Row 22 (source row 44)
$ws.Range("D22").Value = 44446
$ws.Range("J22").Value = 150
$ws.Range("L22").Value = 24000
$ws.Range("M22").Value = 22667
$ws.Range("N22").Value = '$/malla 15 kilos'
$ws.Range("P22").Value = 1511
$ws.Range("Q22").Value = 15

# Row 23 (source row 6)
$ws.Range("D23").Value = 45085
$ws.Range("J23").Value = 200
$ws.Range("L23").Value = 18000
$ws.Range("M23").Value = 17200
$ws.Range("P23").Value = 1147

# Row 24 (source row 57)
$ws.Range("D24").Value = 45112
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 17000
$ws.Range("L24").Value = 18000
$ws.Range("M24").Value = 17400
$ws.Range("P24").Value = 1160

# Row 25 (source row 58)
$ws.Range("D25").Value = 45112
$ws.Range("J25").Value = 380
$ws.Range("K25").Value = 19000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 19395
$ws.Range("P25").Value = 1293

# Row 27 (source row 10)
$ws.Range("D27").Value = 45163
$ws.Range("J27").Value = 350
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = 18686
$ws.Range("P27").Value = 1246

# Row 28 (source row 46)
$ws.Range("D28").Value = 45084
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 18000
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = 18800
$ws.Range("P28").Value = 1253

# Row 29 (source row 31)
$ws.Range("D29").Value = 45128
$ws.Range("J29").Value = 480
$ws.Range("M29").Value = 17521
$ws.Range("P29").Value = 1168

# Row 30 (source row 45)
$ws.Range("D30").Value = 45154
$ws.Range("J30").Value = 400
$ws.Range("K30").Value = 16000
$ws.Range("L30").Value = 18000
$ws.Range("M30").Value = 17150
$ws.Range("P30").Value = 1143

# Row 31 (source row 28)
$ws.Range("D31").Value = 44803
$ws.Range("J31").Value = 400
$ws.Range("K31").Value = 16000
$ws.Range("M31").Value = 16850
$ws.Range("P31").Value = 1123

# Row 32 (source row 17)
$ws.Range("D32").Value = 44761
$ws.Range("J32").Value = 200
$ws.Range("M32").Value = 17400
$ws.Range("P32").Value = 1160

# Row 33 (source row 74)
$ws.Range("D33").Value = 45176
$ws.Range("J33").Value = 290
$ws.Range("K33").Value = 18000
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = 19103
$ws.Range("P33").Value = 1274

# Row 34 (source row 109)
$ws.Range("D34").Value = 45083
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = 18800
$ws.Range("P34").Value = 1253

# Row 35 (source row 65)
$ws.Range("D35").Value = 45049
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 18000
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = 18850
$ws.Range("P35").Value = 1257

# Row 36 (source row 108)
$ws.Range("D36").Value = 44453
$ws.Range("J36").Value = 280
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 22000
$ws.Range("M36").Value = 21286
$ws.Range("P36").Value = 1419

# Row 37 (source row 97)
$ws.Range("D37").Value = 45082
$ws.Range("J37").Value = 400
$ws.Range("K37").Value = 18000
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = 18850
$ws.Range("P37").Value = 1257

# Row 38 (source row 64)
$ws.Range("D38").Value = 45111
$ws.Range("J38").Value = 380
$ws.Range("K38").Value = 19000
$ws.Range("M38").Value = 19395
$ws.Range("P38").Value = 1293

# Row 39 (source row 47)
$ws.Range("D39").Value = 44818
$ws.Range("J39").Value = 230
$ws.Range("K39").Value = 12000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = 12000
$ws.Range("P39").Value = 800

# Row 40 (source row 48)
$ws.Range("D40").Value = 44818
$ws.Range("I40").Value = 'Segunda'
$ws.Range("J40").Value = 170
$ws.Range("K40").Value = 14000
$ws.Range("L40").Value = 14000
$ws.Range("M40").Value = 14000
$ws.Range("P40").Value = 933

# Row 41 (source row 93)
$ws.Range("D41").Value = 45043
$ws.Range("J41").Value = 220
$ws.Range("K41").Value = 18000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = 18909
$ws.Range("P41").Value = 1261

# Row 42 (source row 90)
$ws.Range("D42").Value = 44791
$ws.Range("J42").Value = 300
$ws.Range("M42").Value = 17133
$ws.Range("P42").Value = 1142

# Row 43 (source row 71)
$ws.Range("D43").Value = 44777
$ws.Range("J43").Value = 250
$ws.Range("K43").Value = 14000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = 14400
$ws.Range("P43").Value = 960

# Row 44 (source row 51)
$ws.Range("D44").Value = 44748
$ws.Range("J44").Value = 200
$ws.Range("K44").Value = 16000
$ws.Range("L44").Value = 17000
$ws.Range("M44").Value = 16400
$ws.Range("P44").Value = 1093

# Row 45 (source row 113)
$ws.Range("D45").Value = 44771
$ws.Range("J45").Value = 180
$ws.Range("K45").Value = 18000
$ws.Range("L45").Value = 20000
$ws.Range("M45").Value = 18889
$ws.Range("P45").Value = 1259

# Row 46 (source row 39)
$ws.Range("D46").Value = 45142
$ws.Range("J46").Value = 420
$ws.Range("K46").Value = 15000
$ws.Range("L46").Value = 17000
$ws.Range("M46").Value = 15571
$ws.Range("P46").Value = 1038

# Row 47 (source row 49)
$ws.Range("D47").Value = 44406
$ws.Range("J47").Value = 400
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 22000
$ws.Range("M47").Value = 20850
$ws.Range("P47").Value = 1390

# Row 48 (source row 79)
$ws.Range("D48").Value = 44722
$ws.Range("I48").Value = 'Primera'
$ws.Range("J48").Value = 150
$ws.Range("K48").Value = 18000
$ws.Range("L48").Value = 20000
$ws.Range("M48").Value = 18933
$ws.Range("P48").Value = 1262

# Row 49 (source row 62)
$ws.Range("D49").Value = 45090
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 18000
$ws.Range("L49").Value = 20000
$ws.Range("M49").Value = 19080
$ws.Range("P49").Value = 1272

# Row 50 (source row 105)
$ws.Range("D50").Value = 44749
$ws.Range("J50").Value = 220
$ws.Range("K50").Value = 18000
$ws.Range("L50").Value = 20000
$ws.Range("M50").Value = 19091
$ws.Range("P50").Value = 1273

# Row 51 (source row 92)
$ws.Range("D51").Value = 45134
$ws.Range("J51").Value = 240
$ws.Range("M51").Value = 16500
$ws.Range("P51").Value = 1100

# Row 52 (source row 96)
$ws.Range("D52").Value = 44741
$ws.Range("J52").Value = 250
$ws.Range("M52").Value = 18800
$ws.Range("P52").Value = 1253

# Row 53 (source row 16)
$ws.Range("D53").Value = 45168
$ws.Range("J53").Value = 360
$ws.Range("K53").Value = 17000
$ws.Range("L53").Value = 18000
$ws.Range("M53").Value = 17361
$ws.Range("P53").Value = 1157

# Row 54 (source row 107)
$ws.Range("D54").Value = 44727
$ws.Range("J54").Value = 220
$ws.Range("K54").Value = 16000
$ws.Range("L54").Value = 18000
$ws.Range("M54").Value = 16909
$ws.Range("P54").Value = 1127

# Row 55 (source row 81)
$ws.Range("D55").Value = 44398
$ws.Range("J55").Value = 130
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = 20000
$ws.Range("O55").Value = 'Provincia de Quillota'
$ws.Range("P55").Value = 1333

# Row 56 (source row 95)
$ws.Range("D56").Value = 44775
$ws.Range("J56").Value = 250
$ws.Range("K56").Value = 18000
$ws.Range("L56").Value = 20000
$ws.Range("M56").Value = 19200
$ws.Range("P56").Value = 1280

# Row 57 (source row 116)
$ws.Range("D57").Value = 44790
$ws.Range("J57").Value = 500
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 16000
$ws.Range("M57").Value = 15540
$ws.Range("P57").Value = 1036

# Row 58 (source row 69)
$ws.Range("D58").Value = 44747
$ws.Range("J58").Value = 400
$ws.Range("K58").Value = 17000
$ws.Range("L58").Value = 19000
$ws.Range("M58").Value = 17850
$ws.Range("P58").Value = 1190

# Row 59 (source row 101)
$ws.Range("D59").Value = 44810
$ws.Range("J59").Value = 400
$ws.Range("K59").Value = 17000
$ws.Range("L59").Value = 19000
$ws.Range("M59").Value = 17850
$ws.Range("P59").Value = 1190

# Row 60 (source row 35)
$ws.Range("D60").Value = 45120
$ws.Range("J60").Value = 260
$ws.Range("M60").Value = 17538
$ws.Range("P60").Value = 1169

# Row 61 (source row 50)
$ws.Range("D61").Value = 44755
$ws.Range("J61").Value = 230
$ws.Range("K61").Value = 16000
$ws.Range("M61").Value = 16783
$ws.Range("P61").Value = 1119

# Row 62 (source row 73)
$ws.Range("D62").Value = 44797
$ws.Range("J62").Value = 310
$ws.Range("M62").Value = 18968
$ws.Range("P62").Value = 1265

# Row 63 (source row 11)
$ws.Range("D63").Value = 45126
$ws.Range("J63").Value = 150
$ws.Range("K63").Value = 16000
$ws.Range("L63").Value = 16000
$ws.Range("M63").Value = 16000
$ws.Range("P63").Value = 1067

# Row 64 (source row 12)
$ws.Range("D64").Value = 45126
$ws.Range("J64").Value = 580
$ws.Range("K64").Value = 17000
$ws.Range("L64").Value = 18000
$ws.Range("M64").Value = 17397
$ws.Range("P64").Value = 1160

# Row 65 (source row 76)
$ws.Range("D65").Value = 44736
$ws.Range("J65").Value = 180
$ws.Range("K65").Value = 17000
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = 17889
$ws.Range("P65").Value = 1193

# Row 66 (source row 117)
$ws.Range("D66").Value = 44769
$ws.Range("J66").Value = 400
$ws.Range("K66").Value = 18000
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = 18850
$ws.Range("P66").Value = 1257

# Row 67 (source row 36)
$ws.Range("D67").Value = 44838
$ws.Range("J67").Value = 180
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 16000
$ws.Range("M67").Value = 15444
$ws.Range("P67").Value = 1030

# Row 68 (source row 72)
$ws.Range("D68").Value = 45118
$ws.Range("J68").Value = 400
$ws.Range("L68").Value = 19000
$ws.Range("M68").Value = 17850
$ws.Range("P68").Value = 1190

# Row 69 (source row 70)
$ws.Range("D69").Value = 45114
$ws.Range("J69").Value = 410
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = 17439
$ws.Range("P69").Value = 1163

# Row 70 (source row 42)
$ws.Range("D70").Value = 45086
$ws.Range("J70").Value = 400
$ws.Range("K70").Value = 16000
$ws.Range("M70").Value = 16850
$ws.Range("P70").Value = 1123

# Row 71 (source row 103)
$ws.Range("D71").Value = 44799
$ws.Range("J71").Value = 220
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = 18909
$ws.Range("P71").Value = 1261

# Row 72 (source row 82)
$ws.Range("D72").Value = 45055
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 20000
$ws.Range("L72").Value = 22000
$ws.Range("M72").Value = 20800
$ws.Range("P72").Value = 1387

# Row 73 (source row 104)
$ws.Range("D73").Value = 45093
$ws.Range("J73").Value = 640
$ws.Range("K73").Value = 17000
$ws.Range("L73").Value = 19000
$ws.Range("M73").Value = 17844
$ws.Range("P73").Value = 1190

# Row 74 (source row 86)
$ws.Range("D74").Value = 45121
$ws.Range("J74").Value = 380
$ws.Range("M74").Value = 19053
$ws.Range("P74").Value = 1270

# Row 75 (source row 38)
$ws.Range("D75").Value = 45133
$ws.Range("J75").Value = 310
$ws.Range("K75").Value = 18000
$ws.Range("L75").Value = 20000
$ws.Range("M75").Value = 19032
$ws.Range("P75").Value = 1269

# Row 76 (source row 15)
$ws.Range("D76").Value = 45071
$ws.Range("J76").Value = 290
$ws.Range("K76").Value = 20000
$ws.Range("L76").Value = 22000
$ws.Range("M76").Value = 21172
$ws.Range("P76").Value = 1411

# Row 77 (source row 61)
$ws.Range("D77").Value = 45147
$ws.Range("J77").Value = 410
$ws.Range("K77").Value = 17000
$ws.Range("L77").Value = 18000
$ws.Range("M77").Value = 17366
$ws.Range("P77").Value = 1158

# Row 78 (source row 22)
$ws.Range("D78").Value = 45020
$ws.Range("J78").Value = 240
$ws.Range("K78").Value = 22000
$ws.Range("L78").Value = 23000
$ws.Range("M78").Value = 22625
$ws.Range("N78").Value = '$/malla 17 kilos'
$ws.Range("P78").Value = 1331
$ws.Range("Q78").Value = 17

# Row 79 (source row 91)
$ws.Range("D79").Value = 45092
$ws.Range("J79").Value = 250
$ws.Range("K79").Value = 17000
$ws.Range("L79").Value = 18000
$ws.Range("M79").Value = 17600
$ws.Range("P79").Value = 1173

# Row 80 (source row 25)
$ws.Range("D80").Value = 45070
$ws.Range("J80").Value = 260
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 22000
$ws.Range("M80").Value = 21077
$ws.Range("P80").Value = 1405

# Row 81 (source row 23)
$ws.Range("D81").Value = 44817
$ws.Range("J81").Value = 400
$ws.Range("K81").Value = 16000
$ws.Range("L81").Value = 17000
$ws.Range("M81").Value = 16425
$ws.Range("P81").Value = 1095

# Row 82 (source row 24)
$ws.Range("D82").Value = 44817
$ws.Range("I82").Value = 'Segunda'
$ws.Range("J82").Value = 150
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = 15000
$ws.Range("P82").Value = 1000

# Row 83 (source row 66)
$ws.Range("D83").Value = 44476
$ws.Range("J83").Value = 220
$ws.Range("M83").Value = 20909
$ws.Range("P83").Value = 1394

# Row 84 (source row 83)
$ws.Range("D84").Value = 44365
$ws.Range("J84").Value = 580
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 22000
$ws.Range("M84").Value = 21103
$ws.Range("P84").Value = 1407

# Row 85 (source row 40)
$ws.Range("D85").Value = 44392
$ws.Range("J85").Value = 220
$ws.Range("K85").Value = 23000
$ws.Range("L85").Value = 23000
$ws.Range("M85").Value = 23000
$ws.Range("P85").Value = 1533

# Row 86 (source row 43)
$ws.Range("D86").Value = 45097
$ws.Range("J86").Value = 150
$ws.Range("K86").Value = 16000
$ws.Range("L86").Value = 17000
$ws.Range("M86").Value = 16400
$ws.Range("P86").Value = 1093

# Row 87 (source row 53)
$ws.Range("D87").Value = 45174
$ws.Range("J87").Value = 280
$ws.Range("K87").Value = 18000
$ws.Range("L87").Value = 20000
$ws.Range("M87").Value = 19143
$ws.Range("P87").Value = 1276

# Row 88 (source row 29)
$ws.Range("D88").Value = 45177
$ws.Range("J88").Value = 410
$ws.Range("M88").Value = 17439
$ws.Range("P88").Value = 1163

# Row 89 (source row 32)
$ws.Range("D89").Value = 45146
$ws.Range("J89").Value = 350
$ws.Range("M89").Value = 17343
$ws.Range("P89").Value = 1156

# Row 90 (source row 78)
$ws.Range("D90").Value = 44714
$ws.Range("J90").Value = 200
$ws.Range("L90").Value = 17000
$ws.Range("M90").Value = 16400
$ws.Range("P90").Value = 1093

# Row 91 (source row 102)
$ws.Range("D91").Value = 44391
$ws.Range("J91").Value = 160
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("M91").Value = 20000
$ws.Range("P91").Value = 1333

# Row 92 (source row 9)
$ws.Range("D92").Value = 45169
$ws.Range("J92").Value = 350
$ws.Range("L92").Value = 18000
$ws.Range("M92").Value = 17314
$ws.Range("P92").Value = 1154

# Row 93 (source row 13)
$ws.Range("D93").Value = 44806
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 17000
$ws.Range("L93").Value = 18000
$ws.Range("M93").Value = 17425
$ws.Range("P93").Value = 1162

# Row 94 (source row 27)
$ws.Range("D94").Value = 45091
$ws.Range("J94").Value = 600
$ws.Range("M94").Value = 17583
$ws.Range("P94").Value = 1172

# Row 95 (source row 63)
$ws.Range("D95").Value = 44435
$ws.Range("J95").Value = 140
$ws.Range("K95").Value = 21000
$ws.Range("L95").Value = 23000
$ws.Range("M95").Value = 21714
$ws.Range("P95").Value = 1448

# Row 96 (source row 114)
$ws.Range("D96").Value = 44784
$ws.Range("J96").Value = 220
$ws.Range("K96").Value = 17000
$ws.Range("L96").Value = 18000
$ws.Range("M96").Value = 17455
$ws.Range("P96").Value = 1164

# Row 97 (source row 41)
$ws.Range("D97").Value = 44804
$ws.Range("J97").Value = 310
$ws.Range("K97").Value = 17000
$ws.Range("L97").Value = 18000
$ws.Range("M97").Value = 17581
$ws.Range("P97").Value = 1172

# Row 98 (source row 68)
$ws.Range("D98").Value = 44841
$ws.Range("J98").Value = 580
$ws.Range("K98").Value = 17000
$ws.Range("L98").Value = 18000
$ws.Range("M98").Value = 17448
$ws.Range("P98").Value = 1163

# Row 99 (source row 106)
$ws.Range("D99").Value = 44699
$ws.Range("J99").Value = 150
$ws.Range("K99").Value = 18000
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = 18667
$ws.Range("P99").Value = 1244

# Row 100 (source row 2)
$ws.Range("D100").Value = 45022
$ws.Range("J100").Value = 330
$ws.Range("K100").Value = 20000
$ws.Range("L100").Value = 22000
$ws.Range("M100").Value = 20970
$ws.Range("N100").Value = '$/malla 20 kilos'
$ws.Range("P100").Value = 1048
$ws.Range("Q100").Value = 20

# Row 101 (source row 3)
$ws.Range("D101").Value = 44811
$ws.Range("L101").Value = 18000
$ws.Range("M101").Value = 17425
$ws.Range("P101").Value = 1162

# Row 102 (source row 4)
$ws.Range("D102").Value = 45029
$ws.Range("J102").Value = 220
$ws.Range("K102").Value = 18000
$ws.Range("M102").Value = 18909
$ws.Range("N102").Value = '$/malla 17 kilos'
$ws.Range("P102").Value = 1112
$ws.Range("Q102").Value = 17

# Row 103 (source row 19)
$ws.Range("D103").Value = 44789
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 16000
$ws.Range("M103").Value = 15425
$ws.Range("P103").Value = 1028

# Row 104 (source row 112)
$ws.Range("D104").Value = 45037
$ws.Range("J104").Value = 400
$ws.Range("K104").Value = 20000
$ws.Range("L104").Value = 22000
$ws.Range("M104").Value = 20850
$ws.Range("P104").Value = 1390

# Row 105 (source row 99)
$ws.Range("D105").Value = 45135
$ws.Range("J105").Value = 430
$ws.Range("K105").Value = 16000
$ws.Range("L105").Value = 17000
$ws.Range("M105").Value = 16465
$ws.Range("P105").Value = 1098

# Row 106 (source row 18)
$ws.Range("D106").Value = 44400
$ws.Range("J106").Value = 130
$ws.Range("K106").Value = 24000
$ws.Range("L106").Value = 24000
$ws.Range("M106").Value = 24000
$ws.Range("P106").Value = 1600

# Row 107 (source row 98)
$ws.Range("D107").Value = 44742
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 18000
$ws.Range("L107").Value = 20000
$ws.Range("M107").Value = 18850
$ws.Range("P107").Value = 1257

# Row 108 (source row 30)
$ws.Range("D108").Value = 44798
$ws.Range("J108").Value = 220
$ws.Range("K108").Value = 18000
$ws.Range("L108").Value = 19000
$ws.Range("M108").Value = 18455
$ws.Range("P108").Value = 1230

# Row 109 (source row 87)
$ws.Range("D109").Value = 45127
$ws.Range("K109").Value = 17000
$ws.Range("L109").Value = 18000
$ws.Range("M109").Value = 17400
$ws.Range("P109").Value = 1160

# Row 110 (source row 60)
$ws.Range("D110").Value = 45125
$ws.Range("J110").Value = 580
$ws.Range("K110").Value = 17000
$ws.Range("L110").Value = 18000
$ws.Range("M110").Value = 17397
$ws.Range("P110").Value = 1160

# Row 111 (source row 20)
$ws.Range("D111").Value = 44785
$ws.Range("J111").Value = 400
$ws.Range("K111").Value = 17000
$ws.Range("M111").Value = 17425
$ws.Range("P111").Value = 1162

# Row 112 (source row 33)
$ws.Range("D112").Value = 44399
$ws.Range("J112").Value = 150
$ws.Range("K112").Value = 22000
$ws.Range("M112").Value = 22000
$ws.Range("P112").Value = 1467

# Row 113 (source row 110)
$ws.Range("D113").Value = 45063
$ws.Range("J113").Value = 350
$ws.Range("K113").Value = 23000
$ws.Range("L113").Value = 25000
$ws.Range("M113").Value = 23686
$ws.Range("P113").Value = 1579

# Row 114 (source row 75)
$ws.Range("D114").Value = 44776
$ws.Range("J114").Value = 580
$ws.Range("L114").Value = 19000
$ws.Range("M114").Value = 17897
$ws.Range("P114").Value = 1193

# Row 115 (source row 5)
$ws.Range("D115").Value = 44449
$ws.Range("J115").Value = 220
$ws.Range("K115").Value = 22000
$ws.Range("L115").Value = 24000
$ws.Range("M115").Value = 23091
$ws.Range("P115").Value = 1539

# Row 116 (source row 8)
$ws.Range("D116").Value = 45069
$ws.Range("K116").Value = 18000
$ws.Range("L116").Value = 20000
$ws.Range("M116").Value = 19080
$ws.Range("P116").Value = 1272

# Row 117 (source row 52)
$ws.Range("D117").Value = 45175
$ws.Range("J117").Value = 300
$ws.Range("M117").Value = 19133
$ws.Range("P117").Value = 1276
